# feat: add 2022-Q3 data
#
# 1) Insert a brand-new worksheet "2022-Q3" right after "总计" (this pushes
#    2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 / 2021-Q2 / 2021-Q1 back by one
#    tab position each - no renaming needed, their names stay the same).
# 2) Fill the new sheet with the Q3 fund holdings, matching the layout used
#    by the other quarterly sheets.
# 3) Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q3 and push the previously-existing rows down by one.
# 4) Fix a header typo on the "2021-Q3" sheet: 基金金额 -> 基金规模.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: populate "2022-Q3" with the same layout as the other
# quarterly fund sheets (header row + two fund rows).
# ---------------------------------------------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").Borders.LineStyle = 1
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160

$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A2:A3").Font.Bold = $true
$q3.Range("A2:A3").Borders.LineStyle = 1
$q3.Range("A2:A3").HorizontalAlignment = -4108
$q3.Range("A2:A3").VerticalAlignment = -4160

$q3.Range("B2:B3").NumberFormat = "@"
$q3.Range("D2:G3").NumberFormat = "@"

$q3.Range("B2").Value = "501007"
$q3.Range("C2").Value = "汇添富中证互联网医疗主题指数（LOF）A"
$q3.Range("D2").Value = "0.39"
$q3.Range("E2").Value = "94.69"
$q3.Range("F2").Value = "4.75"
$q3.Range("G2").Value = "0.0185"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "501008"
$q3.Range("C3").Value = "汇添富中证互联网医疗主题指数（LOF）C"
$q3.Range("D3").Value = "0.17"
$q3.Range("E3").Value = "94.69"
$q3.Range("F3").Value = "4.75"
$q3.Range("G3").Value = "0.0081"
$q3.Range("H3").Value = 10

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 3: update the "总计" sheet - push rows 2..7 down to 3..8, then
# write the new 2022-Q3 row on top (row 2). Walking from the bottom up
# avoids clobbering a row before it has been copied down.
# ---------------------------------------------------------------------
for ($r = 7; $r -ge 2; $r--) {
    $dst = $r + 1
    $zongji.Range("B$dst").Value = $zongji.Range("B$r").Value2
    $zongji.Range("C$dst").Value = $zongji.Range("C$r").Value2
    $zongji.Range("D$dst").Value = $zongji.Range("D$r").Value2
}

# row 8 is brand new - give column A the same look (bold/border/centered)
# as the rest of column A by copying the format down from row 7.
$zongji.Range("A7").Copy()
$zongji.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# renumber the index column (0-based, one entry per row)
$zongji.Range("A2").Value = 0
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2
$zongji.Range("A5").Value = 3
$zongji.Range("A6").Value = 4
$zongji.Range("A7").Value = 5
$zongji.Range("A8").Value = 6

# write the new top row for 2022-Q3
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.03

$zongji.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 4: fix the header typo on the "2021-Q3" sheet.
# ---------------------------------------------------------------------
$q3_2021 = $wb.Worksheets.Item("2021-Q3")
$q3_2021.Range("D1").Value = "基金规模"
